$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column C ("max") entirely; this shifts D -> C ("prediction") and E -> D ("rejection-f")
$ws.Range("C1").EntireColumn.Delete() | Out-Null

# Update column B values (row -> value) with the new prediction scores
$bValues = @{
    2  = 28724.4160595666
    3  = 96498.87512547724
    4  = 96484.65499409501
    5  = 81949.7884529935
    6  = 48389.93740341972
    7  = 104549.8678707865
    8  = 102228.8371356339
    9  = 81280.10152616494
    10 = 82870.82047733865
    11 = 87521.31450254522
    12 = 87480.36544053428
    13 = 121875.6974668384
    14 = 50698.79915833386
}

foreach ($row in $bValues.Keys) {
    $ws.Cells.Item($row, 2).Value = $bValues[$row]
}
